$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("A5").Value = 35
$ws1.Range("B5").Value = "FALCON_R"
$ws1.Range("C5").Value = "Capturing reverting ranging markets"
$ws1.Range("D5").Value = "Sell: SuperSmoother crosses Keltner Channels from Top;`nBuy: SupSmooth crosses Keltner Channels from Bottom"
$ws1.Range("E5").Value = "Price goes against Keltner channel with bigger multiplier;`nPrice hit take profit which is either: 3x more than stop loss or we go and cross opposite side of the channel"

$ws1.Rows.Item(5).RowHeight = 90
$ws1.Range("C5:E5").WrapText = $true

$ws1.Activate()
$ws1.Range("E6").Select()

$wb.Save()
